$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Ferrite Bead / Ferrite" BOM row (old row 9). Deleting the
# row shifts every row below it up by one, which already realigns rows
# 10-20 with their new (post-edit) row numbers. ---
$ws.Rows(9).Delete()

# --- Battery/RESET/UART header row: comment text + quantity bumped 3 -> 4 ---
$ws.Range("C4").Value = "Battery, NMI, Reset, UART"
$ws.Range("F4").Value = 4

# --- Cap Semi row: C15 removed from designator list, qty 14 -> 13, and the
# 10nF value dropped from the value list ---
$ws.Range("C6").Value = "C2, C3, C4, C5, C6, C7, C8, C9, C10, C11, C12, C13, C14"
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = ".1uF, 1nF, 2.2uF, 4.7uF, 22pF, 470pF"

# --- Resistor (Res3) row: R23 added to designator list, qty 22 -> 23, and
# the 250k value dropped from the value list ---
$ws.Range("C20").Value = "R1, R2, R3, R4, R5, R6, R7, R8, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19, R20, R21, R22, R23"
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = "1k, 10k, 39, 80.6"

# --- Insert a new BOM row for the SPDT switch (S1) right after the
# resistor row. Copy row 20's formatting down so the new row picks up the
# same borders/font as the rest of the table. ---
$ws.Rows(21).Insert()
$ws.Range("A20:G20").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)
$ws.Range("A21").Value = "SW-SPDT"
$ws.Range("B21").Value = "SPDT Subminiature Toggle Switch, Right Angle Mounting, Vertical Actuation"
$ws.Range("C21").Value = "S1"
$ws.Range("D21").Value = "563-1102-1-ND"
$ws.Range("E21").Value = "SW-SPDT"
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = ""

# --- Crystal oscillator Y1 value: 16Mhz -> 12Mhz ---
$ws.Range("A25").Value = "XTAL 12Mhz"
